$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly entry: insert a row at row 16, pushing existing rows 16-34 down to 17-35.
$ws.Rows.Item(16).Insert()

# Copy the date style (number format) from the row above into the new row's D cell,
# so the new date renders the same way as the rest of the column.
$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 16 with this week's data.
$ws.Range("A16").Value = 5
$ws.Range("B16").Value = "Macroferia Regional de Talca"
$ws.Range("C16").Value = "Maule"
$ws.Range("D16").Value = 44740
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 100112043
$ws.Range("G16").Value = "Pepino dulce"
$ws.Range("H16").Value = "Cultivar IV Región"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 15000
$ws.Range("N16").Value = "$/bandeja 18 kilos"
$ws.Range("O16").Value = "Provincia de Limarí"
$ws.Range("P16").Value = 833
$ws.Range("Q16").Value = 18
$ws.Range("R16").Value = "Hortaliza"
